# Aufgabenliste und Fragen geändert
#
# 1) Add a new worksheet "Fragen" right after "Aufgabenliste" (it becomes the
#    active sheet, matching the source workbook's bookViews/activeTab + the
#    new sheet's tabSelected state).
# 2) Fill in new "Ideen sammeln" mini-table entries (columns I/K, rows 15-19)
#    on the Aufgabenliste sheet.
# 3) Populate the new Fragen sheet with its header + two question rows.
# 4) Restore the per-sheet cell selections recorded in the saved workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Aufgabenliste: new rows in the "Wer? / Erledigt? / Ideen sammeln"
#    side table (columns I/J/K), rows 15-19.
# ---------------------------------------------------------------------
$aufgaben = $wb.Worksheets.Item("Aufgabenliste")

# Row 15: Ulla / 7. /14. fürs Pflichtenheft nachschauen
$aufgaben.Range("I15").Value = "Ulla"
$aufgaben.Range("K7").Copy()
$aufgaben.Range("K15").PasteSpecial(-4122)
$aufgaben.Range("K15").Value = "7. /14. fürs Pflichtenheft nachschauen"

# Row 16: Joana / 8. fürs Pflichtenheft nachschauen
$aufgaben.Range("I16").Value = "Joana"
$aufgaben.Range("K7").Copy()
$aufgaben.Range("K16").PasteSpecial(-4122)
$aufgaben.Range("K16").Value = "8. fürs Pflichtenheft nachschauen"

# Row 17: Bene / Skizze für GUI
$aufgaben.Range("I17").Value = "Bene"
$aufgaben.Range("K7").Copy()
$aufgaben.Range("K17").PasteSpecial(-4122)
$aufgaben.Range("K17").Value = "Skizze für GUI"

# Row 18: Chiara / 11. Systemmodelle nachschauen & 13
$aufgaben.Range("I18").Value = "Chiara"
$aufgaben.Range("K7").Copy()
$aufgaben.Range("K18").PasteSpecial(-4122)
$aufgaben.Range("K18").Value = "11. Systemmodelle nachschauen & 13"

# Row 19: Joana / Programm zum durchlaufen von Programmen
$aufgaben.Range("I19").Value = "Joana"
$aufgaben.Range("K7").Copy()
$aufgaben.Range("K19").PasteSpecial(-4122)
$aufgaben.Range("K19").Value = "Programm zum durchlaufen von Programmen"

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2. Insert the new "Fragen" worksheet right after "Aufgabenliste".
# ---------------------------------------------------------------------
$fragen = $wb.Worksheets.Add($null, $aufgaben)
$fragen.Name = "Fragen"

# Header cell, reuses the same "Erledigt?" header style/text used elsewhere
# in the workbook (e.g. Aufgabenliste!B17 / J5).
$aufgaben.Range("B17").Copy()
$fragen.Range("A1").PasteSpecial(-4122)
$fragen.Range("A1").Value = "Erledigt?"

# Checkbox-style column A down to row 22 (same look as the "Erledigt?"
# column J on Aufgabenliste).
$aufgaben.Range("J6").Copy()
$fragen.Range("A2:A22").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$fragen.Range("B2").Value = "Systemmodelle: UML/Sequenzdiagamme normalerweiße Entwurf?"
$fragen.Range("B3").Value = "Idee für Trace (Antler funktioniert nicht mit Rückgängig machen)"

$fragen.Range("B1").ColumnWidth = 57.83

# ---------------------------------------------------------------------
# 3. Restore selections recorded in the saved file.
# ---------------------------------------------------------------------
$aufgaben.Range("O15").Select()
$fragen.Range("B22").Select()
